$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing trailing rows (50 & 51) with revised values ---
$ws.Range("F50").Value = 62547500

$ws.Range("B51").Value = 209.6799926757812
$ws.Range("C51").Value = 216.8399963378906
$ws.Range("D51").Value = 208.4199981689453
$ws.Range("E51").Value = 215.9499969482422
$ws.Range("F51").Value = 61368300

# --- Append new rows 52-55, cloning the date cell's format (style) from A51 ---
$ws.Range("A51").Copy()
$ws.Range("A52:A55").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A52").Value = 45730
$ws.Range("B52").Value = 213.4900054931641
$ws.Range("C52").Value = 213.9499969482422
$ws.Range("D52").Value = 209.5800018310547
$ws.Range("E52").Value = 211.25
$ws.Range("F52").Value = 60107600

$ws.Range("A53").Value = 45733
$ws.Range("B53").Value = 214
$ws.Range("C53").Value = 215.2200012207031
$ws.Range("D53").Value = 209.9700012207031
$ws.Range("E53").Value = 213.3099975585938
$ws.Range("F53").Value = 48073400

$ws.Range("A54").Value = 45734
$ws.Range("B54").Value = 212.6900024414062
$ws.Range("C54").Value = 215.1499938964844
$ws.Range("D54").Value = 211.4900054931641
$ws.Range("E54").Value = 214.1600036621094
$ws.Range("F54").Value = 42432400

$ws.Range("A55").Value = 45735
$ws.Range("B55").Value = 215.2400054931641
$ws.Range("C55").Value = 218.7599945068359
$ws.Range("D55").Value = 213.75
$ws.Range("E55").Value = 214.2200012207031
$ws.Range("F55").Value = 54336700
